# edit.ps1 - apply "fix bug of OC and USAEx" changes
# Row 7: G7/I7/K7/M7 (US markets closed -> "休市"), O7/P7 (TWD fx rate update)
# Row 8: brand-new trading-day row appended below row 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- Row 7 updates (only the cells that actually changed) ----------
$ws.Range("G7").Value = "休市"
$ws.Range("I7").Value = "休市"
$ws.Range("K7").Value = "休市"
$ws.Range("M7").Value = "休市"
$ws.Range("O7").Value = "'31.416"
$ws.Range("P7").Value = "'-0.013"

# Cells O7/P7 must stay in the sheet's default (unstyled) format; restore it
# after the literal-text assignment above (which otherwise stamps a stray
# quote-prefix style onto the cell).
$ws.Range("Q1").Copy()
$ws.Range("O7:P7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------- Row 8: append the new trading day ----------
$ws.Range("A8").Value = "'2019/05/28"
$ws.Range("B8").Value = "'10312.31"
$ws.Range("C8").Value = "'-21.82"
$ws.Range("D8").Value = "'-0.21 %"
$ws.Range("E8").Value = "1548.68 億"
$ws.Range("F8").Value = "-56.35 億"
$ws.Range("G8").Value = "'0"
$ws.Range("H8").Value = "'0 %"
$ws.Range("I8").Value = "'0"
$ws.Range("J8").Value = "'0 %"
$ws.Range("K8").Value = "'0"
$ws.Range("L8").Value = "'0 %"
$ws.Range("M8").Value = "'0"
$ws.Range("N8").Value = "'0 %"
$ws.Range("O8").Value = "'0"
$ws.Range("P8").Value = "'0"
$ws.Range("Q8").Value = "'41796"
$ws.Range("R8").Value = "'-5195"
$ws.Range("S8").Value = "'1372"
$ws.Range("T8").Value = "'2969"
$ws.Range("U8").Value = "'-204"
$ws.Range("V8").Value = "'2815"
$ws.Range("W8").Value = "'1576"
$ws.Range("X8").Value = "'154"
$ws.Range("Y8").Value = "11.94 / 10.37"
$ws.Range("Z8").Value = "'106.55"
$ws.Range("AA8").Value = "'1.21%"

# Re-apply number/font formatting for row 8 by copying the style from the
# matching column in existing rows, so the new row reuses the workbook's
# existing styles instead of inventing new ones.
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B8:D8").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E8:F8").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("G8:N8").PasteSpecial(-4122)
$ws.Range("Q1").Copy()
$ws.Range("O8:AA8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
